$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header cells: "_old" suffix -> "_FV2404", "_new" suffix -> "_FV2410"
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($c = 1; $c -le 21; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# 2. Turn the used range into an Excel table (ListObject)
$rng = $ws.Range("A1:U69")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = $null

# 3. Freeze the header row (pane split after row 1)
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
